$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '42.369.87'
$ws.Range("E2").Value = '  +2.92%  '
$ws.Range("D3").Value = '2.234.12'
$ws.Range("E3").Value = '  +2.54%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.632'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.32'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.50%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.632'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0939'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.104'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Value = '2.569.86'
$ws.Range("E15").Value = '  +2.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.888'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '2.253.90'
$ws.Range("E18").Value = '  +2.93%  '
$ws.Range("D19").Value = '42.291.84'
$ws.Range("E19").Value = '  +2.90%  '
$ws.Range("D20").Value = '0.0₃0971'
$ws.Range("E20").Value = '  +2.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +2.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.62%  '
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.37'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.80'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.14%  '
$ws.Range("E33").Value = '  +12.32%  '
$ws.Range("E34").Value = '  +4.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0782'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.123'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.67'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0319'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.69%  '
$ws.Range("E41").Value = '  +3.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.200'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.84%  '
$ws.Range("E51").Value = '  +1.33%  '
